$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("summary")
$ws1.Range("C2").Value = 714
$ws1.Range("D2").Value = 689
$ws1.Range("E2").Value = 80.99
$ws1.Range("F2").Value = -1.75
$ws1.Range("G2").Value = 0.1
$ws1.Range("I2").Value = -0.21
$ws1.Range("L2").Value = 1.19
$ws1.Range("C3").Value = 714
$ws1.Range("D3").Value = 689
$ws1.Range("E3").Value = 48.62
$ws1.Range("F3").Value = 0.07
$ws1.Range("H3").Value = 1
$ws1.Range("I3").Value = 0.11
$ws1.Range("J3").Value = 0.27
$ws1.Range("K3").Value = 0.06
$ws1.Range("L3").Value = 1.06
$ws1.Range("C4").Value = 714
$ws1.Range("D4").Value = 669
$ws1.Range("E4").Value = 18.54
$ws1.Range("I4").Value = 0.51
$ws1.Range("L4").Value = 0.82
$ws1.Range("C5").Value = 714
$ws1.Range("D5").Value = 676
$ws1.Range("E5").Value = 63.76
$ws1.Range("F5").Value = -0.69
$ws1.Range("H5").Value = 0.97
$ws1.Range("I5").Value = -0.71
$ws1.Range("J5").Value = 0.32
$ws1.Range("K5").Value = 0.05
$ws1.Range("L5").Value = 1.31
$ws1.Range("C6").Value = 714
$ws1.Range("D6").Value = 643
$ws1.Range("E6").Value = 55.37
$ws1.Range("F6").Value = -0.26
$ws1.Range("I6").Value = 0.81
$ws1.Range("K6").Value = 0.05
$ws1.Range("L6").Value = 0.84
$ws1.Range("C7").Value = 714
$ws1.Range("D7").Value = 623
$ws1.Range("E7").Value = 51.04
$ws1.Range("F7").Value = -0.06
$ws1.Range("J7").Value = 0.35
$ws1.Range("K7").Value = 0.07
$ws1.Range("L7").Value = 1.2
$ws1.Range("C8").Value = 714
$ws1.Range("E8").Value = 47.62
$ws1.Range("F8").Value = 0.12
$ws1.Range("H8").Value = 1
$ws1.Range("I8").Value = -0.01
$ws1.Range("J8").Value = 0.33
$ws1.Range("K8").Value = 0.05
$ws1.Range("L8").Value = 1.01
$ws1.Range("C9").Value = 714
$ws1.Range("E9").Value = 42.47
$ws1.Range("F9").Value = 0.35
$ws1.Range("H9").Value = 1
$ws1.Range("I9").Value = -0.01
$ws1.Range("J9").Value = 0.33
$ws1.Range("K9").Value = 0.03
$ws1.Range("L9").Value = 1.06

$ws2 = $wb.Worksheets.Item("model_fit")
$ws2.Range("B2").Value = 714
$ws2.Range("D2").Value = 6044
$ws2.Range("E2").Value = 6062
$ws2.Range("F2").Value = 6103
$ws2.Range("G2").Value = 0.575
$ws2.Range("H2").Value = 0.408
$ws2.Range("B3").Value = 714
$ws2.Range("D3").Value = 6037
$ws2.Range("E3").Value = 6069
$ws2.Range("F3").Value = 6143
$ws2.Range("G3").Value = 0.581
$ws2.Range("H3").Value = 0.405
